$d = $word.ActiveDocument

# Remove the old email + "(Google talk as well)" text, leaving the
# space run before it ("Email: ") untouched, then insert the new
# email address in its place.
$rng = $d.Content
$rng.Find.Execute("xorcererzc@gmail.com (Google talk as well)", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "", 2)
$rng.InsertAfter("logan.zhou.cn@gmail.com")
